$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shipment record as row 6
$ws.Range("A6").Value = "April 26th, 2022"
$ws.Range("B6").Value = "URI"
$ws.Range("C6").Value = "UPenn"
$ws.Range("D6").Value = "Putnam lab grey large shipper"
$ws.Range("E6").Value = "August 2020 fragments that have been clipped"
$ws.Range("F6").Value = "Emma"
$ws.Range("G6").Value = "Lulu"
$ws.Range("H6").Value = "UPenn shipping empty dry shipper back on May 2nd, 2022"

# Match wrap text styling used in Contents/Notes columns for other rows
$ws.Range("E6").WrapText = $true
$ws.Range("H6").WrapText = $true

# Row auto-sizes slightly taller than the default once wrapped content is present
$ws.Rows.Item(6).RowHeight = 17

$ws.Range("D10").Select() | Out-Null
